$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.803.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.344.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.92%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'240.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.68%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -3.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'72.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.31%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -5.92%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0998"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.90%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'58.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.89%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'32.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.10%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.07%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.693.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.95%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.95%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.904"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.08%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.345.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.78%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.729.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.52%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.74%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.63%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'78.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.52%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'253.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.53%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +8.82%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.01%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -7.05%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'175.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.02%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'22.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.03%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.42%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.50%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.09%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.16%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.61%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.81%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'MultiversX"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'65.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +18.86%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'FTXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'5.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +15.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'9.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.53%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +6.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'18.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.00%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.96%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.64%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'98.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.06%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -4.40%  "
$ws.Range("E51").Style = "Normal"
